$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2074.75
$ws.Range("I4").Value = 1766.3334
$ws.Range("K4").Value = 1766.3334
$ws.Range("M4").Value = -1652.3334
$ws.Range("H33").Value = 154.36363
$ws.Range("I33").Value = 219.5
$ws.Range("J33").Value = 117.14286
$ws.Range("K33").Value = 219.5
$ws.Range("L33").Value = 117.14286
$ws.Range("M33").Value = 9.5
$ws.Range("N33").Value = -575.14286
$ws.Range("H96").Value = 1487.4546
$ws.Range("I96").Value = 1470
$ws.Range("J96").Value = 1508.4
$ws.Range("K96").Value = 4410
$ws.Range("L96").Value = 4525.200000000001
$ws.Range("M96").Value = -3037
$ws.Range("N96").Value = -7271.200000000001
$ws.Range("H112").Value = 2243.6287
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 2250.7942
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 6752.382599999999
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -8968.382599999999
$ws.Range("H132").Value = 1257.6279
$ws.Range("I132").Value = 1159.5
$ws.Range("K132").Value = 3478.5
$ws.Range("M132").Value = -948.5
$ws.Range("H137").Value = 2219.742
$ws.Range("I137").Value = 1660.5714
$ws.Range("J137").Value = 2680.2354
$ws.Range("K137").Value = 4981.7142
$ws.Range("L137").Value = 8040.706200000001
$ws.Range("M137").Value = -2431.7142
$ws.Range("N137").Value = -13140.7062
$ws.Range("H138").Value = 5134.3125
$ws.Range("I138").Value = 6935.857
$ws.Range("J138").Value = 3733.111
$ws.Range("K138").Value = 20807.571
$ws.Range("L138").Value = 11199.333
$ws.Range("M138").Value = -15667.571
$ws.Range("N138").Value = -21479.333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5381.5835
$ws.Range("I32").Value = 4226.282
$ws.Range("K32").Value = 4226.282
$ws.Range("M32").Value = -3939.282
$ws.Range("H61").Value = 1443.1875
$ws.Range("I61").Value = 1422.0769
$ws.Range("J61").Value = 1534.6666
$ws.Range("K61").Value = 1422.0769
$ws.Range("L61").Value = 1534.6666
$ws.Range("M61").Value = -1210.0769
$ws.Range("N61").Value = -1958.6666
$ws.Range("H74").Value = 1720.3334
$ws.Range("I74").Value = 806.7778
$ws.Range("J74").Value = 3090.6667
$ws.Range("K74").Value = 806.7778
$ws.Range("L74").Value = 3090.6667
$ws.Range("M74").Value = 67.22220000000004
$ws.Range("N74").Value = -4838.6667
$ws.Range("H77").Value = 1720.3334
$ws.Range("I77").Value = 806.7778
$ws.Range("J77").Value = 3090.6667
$ws.Range("K77").Value = 4033.889
$ws.Range("L77").Value = 15453.3335
$ws.Range("M77").Value = 334.1110000000003
$ws.Range("N77").Value = -24189.3335
$ws.Range("H130").Value = 14846.167
$ws.Range("J130").Value = 14846.167
$ws.Range("L130").Value = 14846.167
$ws.Range("N130").Value = -24886.167
$ws.Range("H132").Value = 1384.1471
$ws.Range("I132").Value = 1160.88
$ws.Range("J132").Value = 2004.3334
$ws.Range("K132").Value = 3482.64
$ws.Range("L132").Value = 6013.0002
$ws.Range("M132").Value = -952.6400000000003
$ws.Range("N132").Value = -11073.0002
$ws.Range("H136").Value = 1443.1875
$ws.Range("I136").Value = 1422.0769
$ws.Range("J136").Value = 1534.6666
$ws.Range("K136").Value = 4266.2307
$ws.Range("L136").Value = 4603.9998
$ws.Range("M136").Value = -1716.2307
$ws.Range("N136").Value = -9703.9998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 849.1667
$ws.Range("I64").Value = 698.3333
$ws.Range("K64").Value = 698.3333
$ws.Range("M64").Value = -473.3333
$ws.Range("H67").Value = 849.1667
$ws.Range("I67").Value = 698.3333
$ws.Range("K67").Value = 698.3333
$ws.Range("M67").Value = 81.66669999999999
$ws.Range("H107").Value = 514.4286
$ws.Range("I107").Value = 380.2
$ws.Range("J107").Value = 850
$ws.Range("K107").Value = 380.2
$ws.Range("L107").Value = 850
$ws.Range("M107").Value = 1539.8
$ws.Range("N107").Value = -4690
$ws.Range("H134").Value = 10257.571
$ws.Range("I134").Value = 13501.7
$ws.Range("K134").Value = 40505.10000000001
$ws.Range("M134").Value = -37970.10000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2216.3
$ws.Range("I62").Value = 2310.4285
$ws.Range("J62").Value = 1996.6666
$ws.Range("K62").Value = 2310.4285
$ws.Range("L62").Value = 1996.6666
$ws.Range("M62").Value = -1686.4285
$ws.Range("N62").Value = -3244.6666
$ws.Range("H65").Value = 2216.3
$ws.Range("I65").Value = 2310.4285
$ws.Range("J65").Value = 1996.6666
$ws.Range("K65").Value = 11552.1425
$ws.Range("L65").Value = 9983.333000000001
$ws.Range("M65").Value = -8432.1425
$ws.Range("N65").Value = -16223.333
$ws.Range("H132").Value = 3654.8125
$ws.Range("J132").Value = 5594.2
$ws.Range("L132").Value = 16782.6
$ws.Range("N132").Value = -21842.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 458.15384
$ws.Range("I5").Value = 377.9091
$ws.Range("J5").Value = 899.5
$ws.Range("K5").Value = 1133.7273
$ws.Range("L5").Value = 2698.5
$ws.Range("M5").Value = -1021.7273
$ws.Range("N5").Value = -2922.5
$ws.Range("H13").Value = 301
$ws.Range("I13").Value = 301
$ws.Range("K13").Value = 903
$ws.Range("M13").Value = -735
$ws.Range("H131").Value = 17753.586
$ws.Range("J131").Value = 19595.055
$ws.Range("L131").Value = 58785.165
$ws.Range("N131").Value = -68865.16500000001
$ws.Range("H132").Value = 1798.5385
$ws.Range("J132").Value = 2038.1
$ws.Range("L132").Value = 18342.9
$ws.Range("N132").Value = -23402.9
$ws.Range("H135").Value = 458.15384
$ws.Range("I135").Value = 377.9091
$ws.Range("J135").Value = 899.5
$ws.Range("K135").Value = 3401.1819
$ws.Range("L135").Value = 8095.5
$ws.Range("M135").Value = -866.1819
$ws.Range("N135").Value = -13165.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 12557505
$ws.Range("J21").Value = 76673.664
$ws.Range("L21").Value = 76673.664
$ws.Range("N21").Value = -77019.664
$ws.Range("H30").Value = 12557505
$ws.Range("J30").Value = 76673.664
$ws.Range("L30").Value = 76673.664
$ws.Range("N30").Value = -76883.664
$ws.Range("H70").Value = 3509.7273
$ws.Range("I70").Value = 3121.6
$ws.Range("J70").Value = 3833.1667
$ws.Range("K70").Value = 3121.6
$ws.Range("L70").Value = 3833.1667
$ws.Range("M70").Value = -2851.6
$ws.Range("N70").Value = -4373.1667
$ws.Range("H73").Value = 3509.7273
$ws.Range("I73").Value = 3121.6
$ws.Range("J73").Value = 3833.1667
$ws.Range("K73").Value = 3121.6
$ws.Range("L73").Value = 3833.1667
$ws.Range("M73").Value = -2185.6
$ws.Range("N73").Value = -5705.1667
$ws.Range("H126").Value = 39819.11
$ws.Range("I126").Value = 4548.25
$ws.Range("K126").Value = 13644.75
$ws.Range("M126").Value = -11174.75
$ws.Range("H132").Value = 2530.6858
$ws.Range("I132").Value = 1984.909
$ws.Range("J132").Value = 3454.3076
$ws.Range("K132").Value = 5954.727000000001
$ws.Range("L132").Value = 10362.9228
$ws.Range("M132").Value = -3424.727000000001
$ws.Range("N132").Value = -15422.9228
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5440.8
$ws.Range("I7").Value = 4301
$ws.Range("K7").Value = 4301
$ws.Range("M7").Value = -4189
$ws.Range("H55").Value = 398.44827
$ws.Range("I55").Value = 389.2353
$ws.Range("J55").Value = 411.5
$ws.Range("K55").Value = 389.2353
$ws.Range("L55").Value = 411.5
$ws.Range("M55").Value = -216.2353
$ws.Range("N55").Value = -757.5
$ws.Range("H68").Value = 2156.6667
$ws.Range("I68").Value = 1651.1666
$ws.Range("J68").Value = 3167.6667
$ws.Range("K68").Value = 1651.1666
$ws.Range("L68").Value = 3167.6667
$ws.Range("M68").Value = -902.1666
$ws.Range("N68").Value = -4665.6667
$ws.Range("H71").Value = 2156.6667
$ws.Range("I71").Value = 1651.1666
$ws.Range("J71").Value = 3167.6667
$ws.Range("K71").Value = 8255.833000000001
$ws.Range("L71").Value = 15838.3335
$ws.Range("M71").Value = -4511.833000000001
$ws.Range("N71").Value = -23326.3335
$ws.Range("H126").Value = 5440.8
$ws.Range("I126").Value = 4301
$ws.Range("K126").Value = 12903
$ws.Range("M126").Value = -10433
$ws.Range("H132").Value = 6528.5454
$ws.Range("I132").Value = 6913.857
$ws.Range("K132").Value = 20741.571
$ws.Range("M132").Value = -18211.571
$ws.Range("H136").Value = 4419.48
$ws.Range("I136").Value = 3347.7058
$ws.Range("K136").Value = 10043.1174
$ws.Range("M136").Value = -7493.117400000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 59808
$ws.Range("J46").Value = 59808
$ws.Range("L46").Value = 59808
$ws.Range("N46").Value = -60270
$ws.Range("H134").Value = 59808
$ws.Range("J134").Value = 59808
$ws.Range("L134").Value = 179424
$ws.Range("N134").Value = -184494
$ws.Range("H140").Value = 62450
$ws.Range("J140").Value = 62450
$ws.Range("L140").Value = 62450
$ws.Range("N140").Value = -72810
